$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 129
$ws.Range("B129").Value = 0
$ws.Range("C129").Value = 17.20807329
$ws.Range("D129").Value = 14.58552442
$ws.Range("F129").Value = 8.71650363
$ws.Range("G129").Value = 4.67333119
$ws.Range("H129").Value = 5.62085695
$ws.Range("I129").Value = 26.05446455
$ws.Range("J129").Value = 40.99309282
$ws.Range("K129").Value = 0
$ws.Range("L129").Value = 4.78830415
$ws.Range("M129").Value = 6.02733612
$ws.Range("O129").Value = 10.61220495
$ws.Range("P129").Value = 16.34807654
$ws.Range("Q129").Value = 0
$ws.Range("R129").Value = 8.09460938
$ws.Range("S129").Value = 0.08782403
$ws.Range("T129").Value = 1.94960327
$ws.Range("U129").Value = 21.5141204
$ws.Range("V129").Value = 9.30895261
$ws.Range("W129").Value = 5.38962135
$ws.Range("X129").Value = 5.99072799
$ws.Range("Y129").Value = 0
$ws.Range("Z129").Value = 7.50256169
$ws.Range("AA129").Value = 8.54673124
$ws.Range("AB129").Value = 4.98426428
$ws.Range("AD129").Value = 14.02706181
$ws.Range("AE129").Value = 0
$ws.Range("AF129").Value = 10.18305994
$ws.Range("AG129").Value = 44.06894207
$ws.Range("AH129").Value = 8.7611451
$ws.Range("AI129").Value = 18.17558423
$ws.Range("AJ129").Value = 6.89173799
$ws.Range("AK129").Value = 7.11122558
$ws.Range("AL129").Value = 2.56184109
$ws.Range("AM129").Value = 5.37916867
$ws.Range("AN129").Value = 6.04331513
$ws.Range("AO129").Value = 0.87616961
$ws.Range("AP129").Value = 0
$ws.Range("AQ129").Value = 9.589908489999999
$ws.Range("AS129").Value = 4.55979374
$ws.Range("AT129").Value = 13.61055638
$ws.Range("AU129").Value = 2.62965983
$ws.Range("AV129").Value = 3.39230694
$ws.Range("AW129").Value = 6.7438732
$ws.Range("AX129").Value = 11.64725314
$ws.Range("AY129").Value = 8.604010329999999
$ws.Range("BA129").Value = 20.47123762
$ws.Range("BB129").Value = 7.5962332
$ws.Range("BC129").Value = 6.1721437
$ws.Range("BD129").Value = 15.8286717
$ws.Range("BE129").Value = 0

# Row 130
$ws.Range("B130").Value = 0
$ws.Range("C130").Value = 11.40048223
$ws.Range("D130").Value = 21.60637698
$ws.Range("F130").Value = 10.94627302
$ws.Range("G130").Value = 4.00333396
$ws.Range("H130").Value = 10.05645183
$ws.Range("I130").Value = 26.52391811
$ws.Range("J130").Value = 30.45625107
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 5.06222596
$ws.Range("M130").Value = 7.19607686
$ws.Range("O130").Value = 7.28011515
$ws.Range("P130").Value = 16.49292348
$ws.Range("Q130").Value = 3.87629904
$ws.Range("R130").Value = 7.64158813
$ws.Range("S130").Value = 0
$ws.Range("T130").Value = 4.26930726
$ws.Range("U130").Value = 17.12957909
$ws.Range("V130").Value = 5.51843684
$ws.Range("W130").Value = 2.97935059
$ws.Range("X130").Value = 3.29584595
$ws.Range("Y130").Value = 0
$ws.Range("Z130").Value = 7.37237426
$ws.Range("AA130").Value = 5.97360161
$ws.Range("AB130").Value = 8.369071460000001
$ws.Range("AD130").Value = 11.46470357
$ws.Range("AE130").Value = 0
$ws.Range("AF130").Value = 7.15984241
$ws.Range("AG130").Value = 30.79379088
$ws.Range("AH130").Value = 5.06139863
$ws.Range("AI130").Value = 13.96683018
$ws.Range("AJ130").Value = 9.04419399
$ws.Range("AK130").Value = 5.16566561
$ws.Range("AL130").Value = 4.30256142
$ws.Range("AM130").Value = 4.74722383
$ws.Range("AN130").Value = 8.040732520000001
$ws.Range("AO130").Value = 0.03551607
$ws.Range("AP130").Value = 7.12979595
$ws.Range("AQ130").Value = 7.67395586
$ws.Range("AS130").Value = 2.35247588
$ws.Range("AT130").Value = 13.34101353
$ws.Range("AU130").Value = 0
$ws.Range("AV130").Value = 1.05474576
$ws.Range("AW130").Value = 7.06754971
$ws.Range("AX130").Value = 7.02904397
$ws.Range("AY130").Value = 9.73175955
$ws.Range("BA130").Value = 12.42166853
$ws.Range("BB130").Value = 6.61966598
$ws.Range("BC130").Value = 6.15932654
$ws.Range("BD130").Value = 40.37806209
$ws.Range("BE130").Value = 0

# Row 131
$ws.Range("B131").Value = 0
$ws.Range("C131").Value = 6.54229638
$ws.Range("D131").Value = 17.65812753
$ws.Range("F131").Value = 17.24137469
$ws.Range("G131").Value = 4.13257332
$ws.Range("H131").Value = 10.02233859
$ws.Range("I131").Value = 21.03769226
$ws.Range("J131").Value = 21.05491696
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 5.68090728
$ws.Range("M131").Value = 6.86427349
$ws.Range("O131").Value = 4.45401316
$ws.Range("P131").Value = 12.32744833
$ws.Range("Q131").Value = 2.60470253
$ws.Range("R131").Value = 5.68343809
$ws.Range("S131").Value = 0.67238996
$ws.Range("T131").Value = 1.91832383
$ws.Range("U131").Value = 13.01139408
$ws.Range("V131").Value = 2.45594408
$ws.Range("W131").Value = 5.35728562
$ws.Range("X131").Value = 4.51329764
$ws.Range("Y131").Value = 0
$ws.Range("Z131").Value = 6.15100094
$ws.Range("AA131").Value = 9.78661728
$ws.Range("AB131").Value = 5.47142634
$ws.Range("AD131").Value = 8.979015240000001
$ws.Range("AE131").Value = 0
$ws.Range("AF131").Value = 8.93359529
$ws.Range("AG131").Value = 19.34413132
$ws.Range("AH131").Value = 2.07348453
$ws.Range("AI131").Value = 10.10597822
$ws.Range("AJ131").Value = 9.93192825
$ws.Range("AK131").Value = 9.798531029999999
$ws.Range("AL131").Value = 10.05228948
$ws.Range("AM131").Value = 3.44318693
$ws.Range("AN131").Value = 8.922040369999999
$ws.Range("AO131").Value = 0
$ws.Range("AP131").Value = 6.12126127
$ws.Range("AQ131").Value = 7.87571496
$ws.Range("AS131").Value = 10.42740654
$ws.Range("AT131").Value = 15.5100335
$ws.Range("AU131").Value = 0
$ws.Range("AV131").Value = 0
$ws.Range("AW131").Value = 8.01861501
$ws.Range("AX131").Value = 3.27467119
$ws.Range("AY131").Value = 8.113857210000001
$ws.Range("BA131").Value = 5.84931883
$ws.Range("BB131").Value = 7.23781876
$ws.Range("BC131").Value = 4.1368153
$ws.Range("BD131").Value = 32.83434199
$ws.Range("BE131").Value = 0

# Row 132
$ws.Range("B132").Value = 68.98079386000001
$ws.Range("C132").Value = 8.58255484
$ws.Range("D132").Value = 24.70437325
$ws.Range("F132").Value = 18.91333897
$ws.Range("G132").Value = 3.67249539
$ws.Range("H132").Value = 7.87474352
$ws.Range("I132").Value = 24.73973704
$ws.Range("J132").Value = 12.98595827
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 9.06133945
$ws.Range("M132").Value = 6.80630719
$ws.Range("O132").Value = 12.57714401
$ws.Range("P132").Value = 8.62181552
$ws.Range("Q132").Value = 1.61508309
$ws.Range("R132").Value = 7.61090445
$ws.Range("S132").Value = 0.07743877
$ws.Range("T132").Value = 0.16252944
$ws.Range("U132").Value = 9.26529298
$ws.Range("V132").Value = 13.26510061
$ws.Range("W132").Value = 9.42641334
$ws.Range("X132").Value = 2.46608269
$ws.Range("Y132").Value = 0
$ws.Range("Z132").Value = 6.36154568
$ws.Range("AA132").Value = 11.42156943
$ws.Range("AB132").Value = 5.29601895
$ws.Range("AD132").Value = 24.79842656
$ws.Range("AE132").Value = 0
$ws.Range("AF132").Value = 10.28806075
$ws.Range("AG132").Value = 9.868817290000001
$ws.Range("AH132").Value = 0
$ws.Range("AI132").Value = 24.51862952
$ws.Range("AJ132").Value = 10.83665062
$ws.Range("AK132").Value = 7.47412598
$ws.Range("AL132").Value = 11.5933864
$ws.Range("AM132").Value = 4.3115808
$ws.Range("AN132").Value = 6.67586756
$ws.Range("AO132").Value = 0
$ws.Range("AP132").Value = 5.05431991
$ws.Range("AQ132").Value = 8.203993560000001
$ws.Range("AS132").Value = 7.93934126
$ws.Range("AT132").Value = 20.28322983
$ws.Range("AU132").Value = 0
$ws.Range("AV132").Value = 2.03998797
$ws.Range("AW132").Value = 8.34149098
$ws.Range("AX132").Value = 0.37422763
$ws.Range("AY132").Value = 10.00824699
$ws.Range("BA132").Value = 0.77548398
$ws.Range("BB132").Value = 5.0225499
$ws.Range("BC132").Value = 8.798180009999999
$ws.Range("BD132").Value = 25.59059791
$ws.Range("BE132").Value = 0

# Row 133
$ws.Range("B133").Value = 60.83197031
$ws.Range("C133").Value = 7.97757033
$ws.Range("D133").Value = 25.20856944
$ws.Range("F133").Value = 18.38425545
$ws.Range("G133").Value = 4.58487107
$ws.Range("H133").Value = 14.1474077
$ws.Range("I133").Value = 18.97250376
$ws.Range("J133").Value = 6.35579651
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 8.757541870000001
$ws.Range("M133").Value = 6.85349555
$ws.Range("O133").Value = 9.526485510000001
$ws.Range("P133").Value = 5.43345061
$ws.Range("Q133").Value = 0.84431749
$ws.Range("R133").Value = 7.85883252
$ws.Range("S133").Value = 8.60817756
$ws.Range("T133").Value = 0
$ws.Range("U133").Value = 5.97854034
$ws.Range("V133").Value = 21.70381958
$ws.Range("W133").Value = 11.36447413
$ws.Range("X133").Value = 2.5672452
$ws.Range("Y133").Value = 0
$ws.Range("Z133").Value = 5.5208287
$ws.Range("AA133").Value = 10.84520587
$ws.Range("AB133").Value = 7.59345689
$ws.Range("AD133").Value = 26.68647647
$ws.Range("AE133").Value = 0
$ws.Range("AF133").Value = 9.81664662
$ws.Range("AG133").Value = 2.40173939
$ws.Range("AH133").Value = 0
$ws.Range("AI133").Value = 19.36076939
$ws.Range("AJ133").Value = 8.22267759
$ws.Range("AK133").Value = 5.32833836
$ws.Range("AL133").Value = 8.903849900000001
$ws.Range("AM133").Value = 5.39199688
$ws.Range("AN133").Value = 4.66531372
$ws.Range("AO133").Value = 2.91220572
$ws.Range("AP133").Value = 9.607176040000001
$ws.Range("AQ133").Value = 7.64649576
$ws.Range("AS133").Value = 5.71677675
$ws.Range("AT133").Value = 18.94377602
$ws.Range("AU133").Value = 0
$ws.Range("AV133").Value = 6.7439089
$ws.Range("AW133").Value = 9.862051210000001
$ws.Range("AX133").Value = 6.32825181
$ws.Range("AY133").Value = 12.32978222
$ws.Range("BA133").Value = 0
$ws.Range("BB133").Value = 3.13367456
$ws.Range("BC133").Value = 8.892672040000001
$ws.Range("BD133").Value = 18.88431091
$ws.Range("BE133").Value = 0

# Row 134
$ws.Range("A134").Value = "12 06 2020"

# Row 135
$ws.Range("A135").Value = "13 06 2020"

# Row 136
$ws.Range("A136").Value = "14 06 2020"

# Row 137
$ws.Range("A137").Value = "15 06 2020"

# Row 138
$ws.Range("A138").Value = "16 06 2020"
